$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 885.38464
$ws.Range("I9").Value = 278.55554
$ws.Range("J9").Value = 2250.75
$ws.Range("K9").Value = 278.55554
$ws.Range("L9").Value = 2250.75
$ws.Range("M9").Value = -109.55554
$ws.Range("N9").Value = -2588.75

$ws.Range("H33").Value = 13890405
$ws.Range("J33").Value = 2438.3333
$ws.Range("L33").Value = 2438.3333
$ws.Range("N33").Value = -2896.3333

$ws.Range("H98").Value = 3450.8462
$ws.Range("I98").Value = 3185.7896
$ws.Range("K98").Value = 3185.7896
$ws.Range("M98").Value = -1687.7896

$ws.Range("H101").Value = 694.375
$ws.Range("I101").Value = 436.42856
$ws.Range("J101").Value = 2500
$ws.Range("K101").Value = 1309.28568
$ws.Range("L101").Value = 7500
$ws.Range("M101").Value = 312.71432
$ws.Range("N101").Value = -10744

$ws.Range("H106").Value = 3034.3125
$ws.Range("I106").Value = 2749.923
$ws.Range("K106").Value = 2749.923
$ws.Range("M106").Value = -2118.923

$ws.Range("H112").Value = 9602.200000000001
$ws.Range("J112").Value = 4087.5
$ws.Range("L112").Value = 12262.5
$ws.Range("N112").Value = -14478.5

$ws.Range("H122").Value = 3450.8462
$ws.Range("I122").Value = 3185.7896
$ws.Range("K122").Value = 9557.3688
$ws.Range("M122").Value = -7107.3688

$ws.Range("H132").Value = 10634.25
$ws.Range("I132").Value = 7813.931
$ws.Range("K132").Value = 23441.793
$ws.Range("M132").Value = -20911.793

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 1235.5
$ws.Range("I15").Value = 480.66666
$ws.Range("J15").Value = 3500
$ws.Range("K15").Value = 480.66666
$ws.Range("L15").Value = 3500
$ws.Range("M15").Value = -253.66666
$ws.Range("N15").Value = -3954

$ws.Range("H105").Value = 1759.8182
$ws.Range("I105").Value = 1759.8182
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1759.8182
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -12.81819999999993
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 134330
$ws.Range("I134").Value = 250738.75
$ws.Range("K134").Value = 752216.25
$ws.Range("M134").Value = -749681.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7058.8647
$ws.Range("I31").Value = 1188.75
$ws.Range("J31").Value = 13964.883
$ws.Range("K31").Value = 1188.75
$ws.Range("L31").Value = 13964.883
$ws.Range("M31").Value = -893.75
$ws.Range("N31").Value = -14554.883

$ws.Range("H34").Value = 7058.8647
$ws.Range("I34").Value = 1188.75
$ws.Range("J34").Value = 13964.883
$ws.Range("K34").Value = 1188.75
$ws.Range("L34").Value = 13964.883
$ws.Range("M34").Value = -986.75
$ws.Range("N34").Value = -14368.883

$ws.Range("H51").Value = 28642
$ws.Range("J51").Value = 28499.2
$ws.Range("L51").Value = 28499.2
$ws.Range("N51").Value = -29971.2

$ws.Range("H61").Value = 28642
$ws.Range("J61").Value = 28499.2
$ws.Range("L61").Value = 28499.2
$ws.Range("N61").Value = -29195.2

$ws.Range("H94").Value = 617.8461
$ws.Range("I94").Value = 473.25
$ws.Range("K94").Value = 473.25
$ws.Range("M94").Value = -22.25

$ws.Range("H102").Value = 20603.111
$ws.Range("I102").Value = 20219
$ws.Range("K102").Value = 20219
$ws.Range("M102").Value = -17785

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1452.97
$ws.Range("J131").Value = 1471.8541
$ws.Range("L131").Value = 4415.5623
$ws.Range("N131").Value = -14495.5623

$ws.Range("H132").Value = 1895517
$ws.Range("I132").Value = 1660.1428
$ws.Range("J132").Value = 15152515
$ws.Range("K132").Value = 14941.2852
$ws.Range("L132").Value = 136372635
$ws.Range("M132").Value = -12411.2852
$ws.Range("N132").Value = -136377695

$ws.Range("H140").Value = 5905.6
$ws.Range("I140").Value = 3176.6667
$ws.Range("K140").Value = 9530.000100000001
$ws.Range("M140").Value = -4350.000100000001

$ws.Range("H141").Value = 765
$ws.Range("I141").Value = 765
$ws.Range("K141").Value = 2295
$ws.Range("M141").Value = 2885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 531.6
$ws.Range("I9").Value = 350
$ws.Range("J9").Value = 577
$ws.Range("K9").Value = 350
$ws.Range("L9").Value = 577
$ws.Range("M9").Value = -180
$ws.Range("N9").Value = -917

$ws.Range("H12").Value = 2004
$ws.Range("I12").Value = 2004
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2004
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1864
$ws.Range("N12").ClearContents()

$ws.Range("H13").Value = 648.0909
$ws.Range("J13").Value = 871.75
$ws.Range("L13").Value = 871.75
$ws.Range("N13").Value = -1149.75

$ws.Range("H14").Value = 128433.125
$ws.Range("J14").Value = 5775
$ws.Range("L14").Value = 5775
$ws.Range("N14").Value = -6111

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1107.4
$ws.Range("I16").Value = 1147.8462
$ws.Range("J16").Value = 844.5
$ws.Range("K16").Value = 1147.8462
$ws.Range("L16").Value = 844.5
$ws.Range("M16").Value = -977.8462
$ws.Range("N16").Value = -1184.5

$ws.Range("H22").Value = 702.7143
$ws.Range("I22").Value = 669.8333
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 669.8333
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -374.8333
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 702.7143
$ws.Range("I27").Value = 669.8333
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 669.8333
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -562.8333
$ws.Range("N27").Value = -1114

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30200.666
$ws.Range("J80").Value = 30200.666
$ws.Range("L80").Value = 30200.666
$ws.Range("N80").Value = -32196.666

$ws.Range("H83").Value = 30200.666
$ws.Range("J83").Value = 30200.666
$ws.Range("L83").Value = 90601.99800000001
$ws.Range("N83").Value = -100585.998

$ws.Range("H133").Value = 54435.6
$ws.Range("J133").Value = 54435.6
$ws.Range("L133").Value = 54435.6
$ws.Range("N133").Value = -64555.6
